$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.478324770927429
$ws.Range("B1").Value = 1.616544485092163
$ws.Range("C1").Value = 1.726461291313171
$ws.Range("D1").Value = 1.46528422832489
$ws.Range("E1").Value = 1.252141237258911
